$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value2 = $text
    $rng.Style = "Normal"
}

# ---- Row 2 updates (existing row, values change) ----
Set-TextCell "A2" "1326944"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1326944"
$ws.Range("C2").Value = "Digital Marketing Executive"
$ws.Range("D2").Value = "Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "11 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "TIE innovated Solution"

# ---- Row 3 (new) ----
Set-TextCell "A3" "1324636"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1324636"
$ws.Range("C3").Value = "[Impact Fortaleza] -Cost & Quality Planning"
$ws.Range("D3").Value = "Castanhal, Pará, Brasil"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "11 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Petruz Fruity"

# ---- Row 4 (new) ----
Set-TextCell "A4" "1322455"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1322455"
$ws.Range("C4").Value = "[Impact Fortaleza]- Chemical Engineering"
$ws.Range("D4").Value = "Castanhal, PA, Brasil"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "17 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "Petruz Fruity"

# ---- Row 5 (new) ----
Set-TextCell "A5" "1315651"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1315651"
$ws.Range("C5").Value = "Portuguese Language Specialist || Marketing Sector (Flexible RE dates)"
$ws.Range("D5").Value = "Athens, Greece"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "11 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "Travelmyth"

# ---- Column width updates ----
$ws.Columns.Item(3).ColumnWidth = 72.16666666666667
$ws.Columns.Item(4).ColumnWidth = 56.166666666666664
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 24.166666666666668
